$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-17 Sunday" "2024-11-18 Monday"

Replace-Text "77×50=" "84×24="
Replace-Text "37×16=" "41×75="
Replace-Text "30×23=" "59×69="
Replace-Text "89×96=" "18×47="
Replace-Text "55×17=" "12×57="

Replace-Text "62×41=" "53×34="
Replace-Text "89×83=" "17×25="
Replace-Text "26×39=" "38×55="
Replace-Text "30×12=" "12×20="
Replace-Text "12×71=" "46×30="

Replace-Text "20×89=" "45×20="
Replace-Text "25×60=" "28×33="
Replace-Text "43×20=" "28×61="
Replace-Text "13×94=" "25×81="
Replace-Text "97×49=" "45×47="

Replace-Text "41×70=" "50×76="
Replace-Text "81×99=" "27×16="
Replace-Text "14×26=" "29×94="
Replace-Text "21×45=" "70×22="
Replace-Text "40×97=" "29×54="

Replace-Text "35×52=" "27×33="
Replace-Text "47×85=" "88×25="
Replace-Text "79×55=" "29×79="
Replace-Text "29×42=" "61×93="
Replace-Text "91×32=" "39×82="
